$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: updated prices / 1h volume deltas,
# plus a couple of rows that swapped rank order (coin name/link/price/volume together).

$ws.Range("D2").Value = "51.574.68"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "2.743.01"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.17"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.61"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "3.170.16"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "2.741.53"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "51.526.50"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.83"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.05%  "
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.86"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.140"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.35"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0824"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.68"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0349"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +9.84%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +13.22%  "
$ws.Range("D47").Value = "2.100.30"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  -0.62%  "

Write-Host "Applied 105 cell updates"
